$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S3").Value = 46
$ws.Range("S4").Value = 44
$ws.Range("S11").Value = 48
$ws.Range("S14").Value = 47

$ws.Range("S14").Select()
